$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5 so existing rows 5-17 shift down to 6-18.
$ws.Rows.Item(5).Insert()

# Populate the newly-inserted row 5 with the REGISTER entry.
$ws.Range("A5").Value = "REGISTER"
$ws.Range("B5").Value = "sword"
$ws.Range("C5").Value = "spear"
$ws.Range("D5").Value = "END"

# Update the selected cell to match the authored state.
$ws.Range("D5").Select()
